# Applies the edits described in the commit diff to "Requisitos de la
# aplicación web.docx":
#   1. Fill in the first blank ("____" -> "ropa") and fix the agreement
#      of "clasificarlos" -> "clasificarla" in the intro paragraph.
#   2. Touch the "Usaremos APIs externas..." paragraph (no visible text
#      change) so the run layout collapses the way Word leaves it after
#      re-editing the sentence (the spell-check proofErr markers go away).
#   3. Replace the three blank "En primer lugar / Además / Finalmente"
#      placeholder paragraphs with the seven fully written paragraphs
#      describing the functional requirements of each page.

$d = $word.ActiveDocument

function Replace-ParagraphText($paragraph, [string]$newText) {
    # Replace the visible contents of a paragraph (everything except its
    # trailing paragraph mark) using Find/Replace so that runs spanning
    # the whole match collapse into a single, uniformly formatted run -
    # mirrors what Word does when the sentence is retyped.
    $start = $paragraph.Range.Start
    $end = $paragraph.Range.End
    $body = $d.Range($start, $end - 1)
    $oldText = $body.Text
    [void]$body.Find.Execute($oldText, $false, $false, $false, $false, `
        $false, $true, 1, $false, $newText, 2)
}

# ---------------------------------------------------------------------
# 1) Intro paragraph: fill in the product blank and fix agreement.
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs(2)
$introRange = $introPara.Range
[void]$introRange.Find.Execute("____", $true, $false, $false, $false, `
    $false, $true, 1, $false, "ropa", 2)

$introRange2 = $introPara.Range
[void]$introRange2.Find.Execute("clasificarlos", $true, $false, $false, `
    $false, $false, $true, 1, $false, "clasificarla", 2)

# ---------------------------------------------------------------------
# 2) "Usaremos APIs externas..." paragraph: re-assert the same text so
#    the split proof-reading runs collapse into a single run.
# ---------------------------------------------------------------------
$apiPara = $d.Paragraphs(9)
Replace-ParagraphText $apiPara "Usaremos APIs externas para poder aceptar pagos, más concretamente la API de PayPal."

# ---------------------------------------------------------------------
# 3) Replace the three placeholder paragraphs with the seven completed
#    ones describing the functional requirements.
# ---------------------------------------------------------------------
# The "Además,___" and "Finalmente,___" placeholders both start with a
# stray <w:proofErr/> marker that sits *before* their first run, so
# rewriting their text in place would leave the marker behind. Drop both
# paragraphs outright (from the bottom up, so indices stay valid) and
# rebuild everything after "En primer lugar,___" from scratch instead.
$finalmentePara = $d.Paragraphs(23)
$finalmentePara.Range.Delete()
$ademasPara = $d.Paragraphs(22)
$ademasPara.Range.Delete()

$p1 = $d.Paragraphs(21)
Replace-ParagraphText $p1 "En primer lugar, existirá una barra de navegación en la parte superior de la página, que permitirá al usuario navegar para encontrar las diferentes funcionalidades. En esta barra se podrán encontrar el acceso a los usuarios, un buscador, las diferentes categorías de ropa, un enlace a la página principal y otro al carrito de la compra. A continuación, se visualizará un hero que muestre enlaces a diferentes categorías destacadas y cambie automáticamente de item cada cierto tiempo. Se mostrará una fila de iconos que informen al usuario sobre los envíos, política de devolución y datos de interés. Seguidamente, se encontrará el grid principal con ítems destacados. Al final de la página se encontrará un footer que recoja información de interés"

# Insert six new empty paragraphs right after paragraph 1, in order, to
# host the remaining six paragraphs (search-results page, product detail
# page, cart page, checkout page, login/registration page, and the
# closing "Finalmente" paragraph about the user's profile page).
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()

$d.Paragraphs(22).Range.Text = "La página de búsqueda contará con un grid principal que muestre los resultados. Además, se mostrarán diferentes opciones de filtrado en un menú lateral."
$d.Paragraphs(23).Range.Text = "Al seleccionar un elemento del grid, se mostrará una página de detalles que informará al usuario de la información del producto, además de permitir añadirlo al carrito de la compra. Debajo de la información del producto se situará un grid con productos relacionados."
$d.Paragraphs(24).Range.Text = "La página del carrito de la compra mostrará los ítems añadidos en forma de tabla y permitirá eliminar elementos. También habrá un botón para procesar el pedido."
$d.Paragraphs(25).Range.Text = "A la hora de procesar el pedido, se mostrará un formulario con la información necesaria para tramitar el pedido, además de la selección de pago."
$d.Paragraphs(26).Range.Text = "La página de inicio de sesión/registro mostrará un formulario con los datos necesarios para el inicio de sesión/registro."
$d.Paragraphs(27).Range.Text = "Finalmente, la página del usuario mostrará un formulario editable con la información del usuario."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
